# Resume edit (commit: "100924 commit - updated resume once more")
#
# The sentence:
#   "...attention to detail during my experience with DoD programs and
#    operational deployments. Have successfully led..."
# becomes:
#   "...attention to detail throughout my experiences within DoD programs
#    and their operational deployments. Have successfully led..."
#
# i.e. "during my experience with DoD programs and operational deployments"
#   -> "throughout my experiences within DoD programs and their operational
#       deployments"
#
# The surrounding text (and its Arimo/20pt/2C2C2C formatting) is untouched;
# only this clause is swapped out, scoped tightly so nothing else in the
# paragraph/document is disturbed.

$d = $word.ActiveDocument

$oldText = "during my experience with DoD programs and operational deployments"
$newText = "throughout my experiences within DoD programs and their operational deployments"

$found = $d.Content.Find.Execute($oldText, $true, $false, $false, $false, $false, $true, 1, $false, $newText, 2)

if (-not $found) {
    throw "Could not locate the target sentence fragment to replace."
}

Write-Output "Replaced: $found"
